$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their original text formatting (e.g. trailing
# zeros like "1.00" or "332.40") instead of being auto-coerced to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "51.146.64"
$ws.Range("E2").Value = "  +2.52%  "
$ws.Range("D3").Value = "2.749.12"
$ws.Range("E3").Value = "  +2.75%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "114.91"
$ws.Range("E5").Value = "  +1.01%  "
$ws.Range("D6").Value = "332.40"
$ws.Range("E6").Value = "  +2.04%  "
$ws.Range("D7").Value = "0.531"
$ws.Range("E7").Value = "  +0.54%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  +3.64%  "
$ws.Range("D10").Value = "41.35"
$ws.Range("E10").Value = "  +0.71%  "
$ws.Range("D11").Value = "20.16"
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("E12").Value = "  +0.72%  "
$ws.Range("E13").Value = "  +3.01%  "
$ws.Range("D14").Value = "7.64"
$ws.Range("E14").Value = "  +4.00%  "
$ws.Range("D15").Value = "3.180.68"
$ws.Range("E15").Value = "  +3.13%  "
$ws.Range("D16").Value = "2.741.88"
$ws.Range("E16").Value = "  +2.47%  "
$ws.Range("D17").Value = "0.887"
$ws.Range("E17").Value = "  +1.42%  "
$ws.Range("D18").Value = "51.179.44"
$ws.Range("E18").Value = "  +2.70%  "
$ws.Range("D19").Value = "13.63"
$ws.Range("E19").Value = "  +3.74%  "
$ws.Range("D20").Value = "3.02"
$ws.Range("D21").Value = "6.85"
$ws.Range("E21").Value = "  +0.93%  "
$ws.Range("D22").Value = "0.0₃0962"
$ws.Range("E22").Value = "  +0.42%  "
$ws.Range("D23").Value = "279.89"
$ws.Range("E23").Value = "  +0.46%  "
$ws.Range("E24").Value = "  -2.71%  "
$ws.Range("D25").Value = "2.62"
$ws.Range("E25").Value = "  +1.77%  "
$ws.Range("D26").Value = "27.00"
$ws.Range("E26").Value = "  +0.51%  "
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("D28").Value = "10.32"
$ws.Range("E28").Value = "  +2.01%  "
$ws.Range("E29").Value = "  -0.77%  "
$ws.Range("E30").Value = "  -1.18%  "
$ws.Range("D31").Value = "35.55"
$ws.Range("E31").Value = "  -1.82%  "
$ws.Range("D32").Value = "49.99"
$ws.Range("E32").Value = "  -0.71%  "
$ws.Range("D33").Value = "5.62"
$ws.Range("E33").Value = "  +2.48%  "
$ws.Range("D34").Value = "0.0825"
$ws.Range("E34").Value = "  +1.46%  "
$ws.Range("D35").Value = "19.37"
$ws.Range("E35").Value = "  -1.12%  "
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "2.11"
$ws.Range("E37").Value = "  +1.48%  "
$ws.Range("D38").Value = "5.04"
$ws.Range("E38").Value = "  -1.48%  "
$ws.Range("D39").Value = "3.22"
$ws.Range("E39").Value = "  +1.56%  "
$ws.Range("D40").Value = "129.32"
$ws.Range("E40").Value = "  +3.92%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "0.0354"
$ws.Range("E41").Value = "  +11.21%  "
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").Value = "23.61"
$ws.Range("E42").Value = "  +3.68%  "
$ws.Range("D43").Value = "2.30"
$ws.Range("E43").Value = "  +3.52%  "
$ws.Range("E44").Value = "  +0.50%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.114.51"
$ws.Range("E45").Value = "  -0.06%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "3.40"
$ws.Range("E46").Value = "  +3.04%  "
$ws.Range("D47").Value = "2.21"
$ws.Range("E47").Value = "  +8.95%  "
$ws.Range("D48").Value = "2.29"
$ws.Range("E48").Value = "  +1.04%  "
$ws.Range("D49").Value = "5.53"
$ws.Range("E49").Value = "  +2.86%  "
$ws.Range("E50").Value = "  +0.69%  "
$ws.Range("D51").Value = "1.54"
$ws.Range("E51").Value = "  +9.36%  "
